$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Devices sheet: refresh UPDATED_TS (column K) on a handful of existing
#    devices that were touched/re-saved during this session, and append a
#    brand-new "SMART_LIGHT" device row (Calex Light / SL001).
# ---------------------------------------------------------------------------
$devices = $wb.Worksheets.Item("Devices")

$devices.Cells.Item(2, 11).Value  = "2025-07-09T00:00:16.110652697+02:00[Europe/Amsterdam]"   # LI001
$devices.Cells.Item(3, 11).Value  = "2025-07-08T21:28:19.600290268+02:00[Europe/Amsterdam]"   # LI002
$devices.Cells.Item(16, 11).Value = "2025-07-08T21:28:19.807904786+02:00[Europe/Amsterdam]"   # LI015
$devices.Cells.Item(36, 11).Value = "2025-07-08T09:15:05.953433389+02:00[Europe/Amsterdam]"   # WA001 / FamilyWM
$devices.Cells.Item(37, 11).Value = "2025-07-09T00:00:15.078798183+02:00[Europe/Amsterdam]"   # DR001

$newRow = 49
$devices.Cells.Item($newRow, 1).Value  = "SMART_LIGHT"
$devices.Cells.Item($newRow, 2).Value  = "SL001"
$devices.Cells.Item($newRow, 3).Value  = "Calex Light"
$devices.Cells.Item($newRow, 4).Value  = "Calex"
$devices.Cells.Item($newRow, 5).Value  = "Calex A60E27"
$devices.Cells.Item($newRow, 6).Value  = $false
$devices.Cells.Item($newRow, 7).Value  = 1050.0
$devices.Cells.Item($newRow, 8).Value  = 1050.0
$devices.Cells.Item($newRow, 9).Value  = "on, off, setMode, status"
$devices.Cells.Item($newRow, 10).Value = "2025-07-08T14:48:46.661059400+02:00[Europe/Amsterdam]"
$devices.Cells.Item($newRow, 11).Value = "2025-07-08T16:13:17.495144517+02:00[Europe/Amsterdam]"
$devices.Cells.Item($newRow, 12).Value = "N/A"
$devices.Cells.Item($newRow, 13).Value = 0.0
$devices.Cells.Item($newRow, 14).Value = 0.0
$devices.Cells.Item($newRow, 15).Value = 0.0
$devices.Cells.Item($newRow, 16).Value = "None"

# ---------------------------------------------------------------------------
# 2) Sense_Control sheet: row 3 (LIGHT LI015 -> LIGHT LIs001 sensor link) was
#    rewritten in this session, losing its explicit formatting in the
#    process.
# ---------------------------------------------------------------------------
$senseControl = $wb.Worksheets.Item("Sense_Control")
$senseControl.Rows.Item(3).ClearFormats()
$senseControl.Cells.Item(3, 1).Value = "LIGHT"
$senseControl.Cells.Item(3, 2).Value = "LI015"
$senseControl.Cells.Item(3, 3).Value = "LIGHT"
$senseControl.Cells.Item(3, 4).Value = "LIs001"
$senseControl.Cells.Item(3, 5).Value = 1400.0
$senseControl.Cells.Item(3, 6).Value = 1400.0

# ---------------------------------------------------------------------------
# 3) New "Smart_Light_Control" sheet (light modes/effects catalogue) appended
#    after the existing "Tasks" sheet.
# ---------------------------------------------------------------------------
$smartLight = $wb.Worksheets.Add()
$smartLight.Name = "Smart_Light_Control"
$tasks = $wb.Worksheets.Item("Tasks")
$smartLight.Move([System.Reflection.Missing]::Value, $tasks)

$smartLight.Cells.Item(1, 1).Value = "MODE_NAME"
$smartLight.Cells.Item(1, 2).Value = "R"
$smartLight.Cells.Item(1, 3).Value = "G"
$smartLight.Cells.Item(1, 4).Value = "B"
$smartLight.Cells.Item(1, 5).Value = "IS_DEFAULT"
$smartLight.Cells.Item(1, 6).Value = "EFFECT_NAME"
$smartLight.Cells.Item(1, 7).Value = "TYPE"
$smartLight.Cells.Item(1, 8).Value = "PARAMS"
